$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell's formatting (bold, centered, bordered) onto A2
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new row of values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1.014133066666666
$ws.Range("C2").Value = 1.741524066666667
$ws.Range("D2").Value = 0.5509200666666665
$ws.Range("E2").Value = 1.531733333333335
$ws.Range("F2").Value = -0.6349063333333298
$ws.Range("G2").Value = -0.2803609999999992
